# Delete the first slide ("Objective: Multi-sense Language Model" /
# Write-Assistant bullet slide). The remaining slide (the Hexagon word-sense
# diagram, formerly slide 2 / sldId 257) becomes the sole slide in the deck.
$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
